$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 1087
$wsExhibit.Range("F5").Value = 3525
$wsExhibit.Range("F10").Value = 19
$wsExhibit.Range("F13").Value = 209
$wsExhibit.Range("F14").Value = 42
$wsExhibit.Range("F16").Value = 2812
$wsExhibit.Range("F17").Value = 1125

# Sheet "全部类型" (all types) - same events mirrored one row further down,
# update "想去人数" (want-to-go count) column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1087
$wsAll.Range("F6").Value = 3526
$wsAll.Range("F12").Value = 19
$wsAll.Range("F15").Value = 209
$wsAll.Range("F16").Value = 42
$wsAll.Range("F18").Value = 2812
$wsAll.Range("F19").Value = 1125
